# Autogenerated on Fri Mar 20 2015 00:16:06 GMT+0000 (Coordinated Universal Time)
#
# The "Enterprises (absolute #)" row and the "Enterprises density (per 1000
# people)" row were swapped, in both the "Source Type: Statistical
# Institution" table (rows 10-11) and the "Source Type: SME Associations"
# table (rows 32-33), so the density row now appears before the
# absolute-count row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking values in this sheet are stored as TEXT (not numbers), so
# force text number-formatting before writing them back, otherwise the COM
# layer auto-coerces strings like "20.1" / "429035" into Doubles.
$ws.Range("D10:D11").NumberFormat = "@"
$ws.Range("B32:D33").NumberFormat = "@"

# --- Table 1 (Source Type: Statistical Institution), rows 10 & 11 ---
$a10 = $ws.Range("A10").Value()
$d10 = $ws.Range("D10").Value()
$a11 = $ws.Range("A11").Value()
$d11 = $ws.Range("D11").Value()

$ws.Range("A10").Value = $a11
$ws.Range("D10").Value = $d11
$ws.Range("A11").Value = $a10
$ws.Range("D11").Value = $d10

# --- Table 2 (Source Type: SME Associations), rows 32 & 33 ---
$a32 = $ws.Range("A32").Value()
$b32 = $ws.Range("B32").Value()
$c32 = $ws.Range("C32").Value()
$d32 = $ws.Range("D32").Value()

$a33 = $ws.Range("A33").Value()
$b33 = $ws.Range("B33").Value()
$c33 = $ws.Range("C33").Value()
$d33 = $ws.Range("D33").Value()

$ws.Range("A32").Value = $a33
$ws.Range("B32").Value = $b33
$ws.Range("C32").Value = $c33
$ws.Range("D32").Value = $d33

$ws.Range("A33").Value = $a32
$ws.Range("B33").Value = $b32
$ws.Range("C33").Value = $c32
$ws.Range("D33").Value = $d32
